$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row 4: IC-C3 (a third "C" reef site) between IC-C2 and IC-U1
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "IC-C3"
$ws.Range("B4").Value = "C"
$ws.Range("C4").Value = 24.53109
$ws.Range("D4").Value = -81.48502

# Insert new row 7: IC-U3 (a third "U" reef site) between IC-U2 and IC-Z1
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "IC-U3"
$ws.Range("B7").Value = "U"
$ws.Range("C7").Value = 24.54047
$ws.Range("D7").Value = -81.44465

# Insert new row 10: IC-Z3 (a third "Z" reef site) between IC-Z2 and the summary rows
$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = "IC-Z3"
$ws.Range("B10").Value = "Z"
$ws.Range("C10").Value = 24.52757
$ws.Range("D10").Value = -81.49843

# Set the page to print in portrait orientation
$ws.PageSetup.Orientation = 1
